# Apply "Constraint 4" data edits to the Data sheet, and update the
# active-sheet/active-cell selection state to match the authored commit.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Explanation")

# --- Data edits: columns D:G ("Constraint 4") ---------------------------
# Rows 3-8 (node-pair group 1,1 .. 1,6): 8,8,8,8 -> 5,6,7,2
foreach ($r in 3..8) {
    $ws1.Cells.Item($r, 4).Value = 5   # D
    $ws1.Cells.Item($r, 5).Value = 6   # E
    $ws1.Cells.Item($r, 6).Value = 7   # F
    $ws1.Cells.Item($r, 7).Value = 2   # G
}

# Rows 9-14 (node-pair group 2,1 .. 2,6): 15,15,15,15 -> 3,4,8,1
foreach ($r in 9..14) {
    $ws1.Cells.Item($r, 4).Value = 3   # D
    $ws1.Cells.Item($r, 5).Value = 4   # E
    $ws1.Cells.Item($r, 6).Value = 8   # F
    $ws1.Cells.Item($r, 7).Value = 1   # G
}

# --- View/selection state -------------------------------------------------
# Explanation was the active tab with U11 selected; after the edit, Data
# becomes the active tab with F9 selected, and Explanation's selection is
# left at U11 (no longer the active tab).
$ws2.Range("U11").Select()
$ws1.Select()
$ws1.Range("F9").Select()
